# edit.ps1 -- applies the "Week 7 -> Week 8 start" content edits described
# by the diff:
#   1. "Kernels I" + "I" (two runs) -> single run "Kernels II"
#   2. Italicize the lone "x" in "Cost: solve function of ...Tx"
#   3. Merge the " " and "LR/SVM w/ no kernel" runs (in the
#      "Small n, large m -> LR/SVM w/ no kernel (slow w/ SVM)" bullet)
#      without touching the trailing " (slow w/ SVM)" run.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "Kernels I" + "I" -> "Kernels II"
#    The paragraph currently holds two runs whose text concatenates to
#    "Kernels II". Find/Replace on that paragraph's range rewrites the
#    first run's text and removes the (now redundant) second run.
# ---------------------------------------------------------------------
$kernelsPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "Kernels II") {
        $kernelsPara = $p
    }
}
if ($kernelsPara -ne $null) {
    $kRange = $kernelsPara.Range
    $kRange.Find.ClearFormatting()
    $kRange.Find.Execute("Kernels I" + "I", $true, $false, $false, $false, `
        $false, $true, 1, $false, "Kernels II", 2) | Out-Null
}

# ---------------------------------------------------------------------
# 2) Italicize the stand-alone "x" run (theta^T x).
# ---------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("Cost: solve function of")) {
        $pr = $p.Range
        $n = $pr.Characters.Count
        for ($i = 1; $i -le $n; $i++) {
            $ch = $pr.Characters($i)
            if ($ch.Text -eq "x") {
                $absStart = $pr.Start + $ch.Start
                $absEnd = $pr.Start + $ch.End
                $xRange = $d.Range($absStart, $absEnd)
                $xRange.Font.Italic = $true
            }
        }
    }
}

# ---------------------------------------------------------------------
# 3) Merge the " " run with the following "LR/SVM w/ no kernel" run,
#    leaving the trailing " (slow w/ SVM)" run untouched.
#
#    This engine coalesces a just-edited run with any immediately
#    following run(s) that share identical formatting, so a direct
#    replace would also swallow " (slow w/ SVM)". To avoid that we
#    briefly give that trailing run different formatting (Bold) so it
#    is not coalesced, perform the merge, then clear the temporary
#    Bold again (a pure formatting change, which does not trigger the
#    run-coalescing pass).
# ---------------------------------------------------------------------
$targetPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("Small n, large m")) {
        $targetPara = $p
    }
}

if ($targetPara -ne $null) {
    $full = $targetPara.Range
    $fullText = $full.Text
    $markerText = " LR/SVM w/ no kernel"
    $tailText = " (slow w/ SVM)"

    $mergeStart = $fullText.IndexOf($markerText)
    if ($mergeStart -ge 0) {
        $mergeEnd = $mergeStart + $markerText.Length
        $tailStart = $mergeEnd
        $tailEnd = $tailStart + $tailText.Length

        # Step 1: shield the trailing run from the coalescing pass.
        $tailAbsStart = $full.Start + $tailStart
        $tailAbsEnd = $full.Start + $tailEnd
        $tailRange = $d.Range($tailAbsStart, $tailAbsEnd)
        if ($tailRange.Text -eq $tailText) {
            $tailRange.Bold = $true

            # Step 2: replace the " " + "LR/SVM w/ no kernel" runs via a
            # placeholder (a real content change, so the engine performs
            # its run-coalescing cleanup), then restore the exact text.
            $mergeAbsStart = $full.Start + $mergeStart
            $mergeAbsEnd = $full.Start + $mergeEnd
            $mergeRange = $d.Range($mergeAbsStart, $mergeAbsEnd)
            $placeholder = "@@MERGE_PLACEHOLDER@@"
            $mergeRange.Text = $placeholder

            $full2 = $targetPara.Range
            $phAbsStart = $full2.Start + $mergeStart
            $phAbsEnd = $phAbsStart + $placeholder.Length
            $phRange = $d.Range($phAbsStart, $phAbsEnd)
            $phRange.Text = $markerText

            # Step 3: remove the temporary Bold shield so the trailing
            # run's formatting matches the original exactly again.
            $full3 = $targetPara.Range
            $tailAbsStart2 = $full3.Start + $mergeStart + $markerText.Length
            $tailAbsEnd2 = $tailAbsStart2 + $tailText.Length
            $tailRange2 = $d.Range($tailAbsStart2, $tailAbsEnd2)
            if ($tailRange2.Text -eq $tailText) {
                $tailRange2.Bold = $false
            }
        }
    }
}

Write-Output "done"
